# issue #5: add legislator_id, name, date into dataframe
# Target sheet: "股票" (Stock), which is the 6th worksheet in this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# ---- Header row (row 1): add three new headers in H1:J1 ----
$ws.Cells.Item(1,8).Value  = "date"
$ws.Cells.Item(1,9).Value  = "legislator_name"
$ws.Cells.Item(1,10).Value = "legislator_id"

# Copy header style (bold + border, same as existing header cells) onto the new header cells
$ws.Cells.Item(1,7).Copy()
$ws.Range($ws.Cells.Item(1,8), $ws.Cells.Item(1,10)).PasteSpecial(-4122)

# ---- Data rows (2 through 11): populate date / legislator_name / legislator_id ----
for ($r = 2; $r -le 11; $r++) {
    # Prefix with a single quote so Excel stores the date as literal text "2012-05-01"
    # rather than auto-converting it to a date serial number.
    $ws.Cells.Item($r,8).Value  = "'2012-05-01"
    $ws.Cells.Item($r,9).Value  = "黃偉哲"
    $ws.Cells.Item($r,10).Value = 1367

    # Re-apply the same (plain) formatting used by the rest of the data row to column H,
    # so it matches the body style of the other cells instead of a date format.
    $ws.Cells.Item($r,7).Copy()
    $ws.Cells.Item($r,8).PasteSpecial(-4122)
}
